# Rename "Sheet1" to "测试" and make it the active/selected sheet
# (moving the active tab away from Sheet2).

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Name = "测试"
$sheet1.Activate()
$sheet1.Select()
